# Update "想去人数" (want-to-go count) figures on the 展览 and 全部类型 sheets
# to reflect the newly scraped values.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 7455
$ws1.Range("F5").Value = 224
$ws1.Range("F6").Value = 1124
$ws1.Range("F8").Value = 17
$ws1.Range("F9").Value = 118
$ws1.Range("F10").Value = 30

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 7455
$ws4.Range("F5").Value = 224
$ws4.Range("F6").Value = 1124
$ws4.Range("F9").Value = 17
$ws4.Range("F10").Value = 118
$ws4.Range("F11").Value = 30
